$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B10 content (Organization website) from "www.stat.kg " to "www.stat.gov.kg"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Update selection to B10
$ws.Range("B10").Select()
